# Auto-generated script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking values in columns D and E are stored as text (matching source data),
# so Excel does not auto-convert them to numbers (losing formatting / precision).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "39.369.21"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "2.156.71"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "228.95"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("D7").Value = "63.04"
$ws.Range("E7").Value = "  +3.91%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("D10").Value = "0.0861"
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "16.10"
$ws.Range("E12").Value = "  +7.47%  "
$ws.Range("D13").Value = "2.475.12"
$ws.Range("E13").Value = "  +3.23%  "
$ws.Range("D14").Value = "22.29"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").Value = "0.819"
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "2.151.33"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "39.431.16"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").Value = "72.18"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("D21").Value = "0.0₃0856"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("D22").Value = "228.40"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "9.79"
$ws.Range("D27").Value = "171.72"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "19.66"
$ws.Range("E29").Value = "  +2.51%  "
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").Value = "2.57"
$ws.Range("E31").Value = "  +9.19%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "4.63"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("D34").Value = "4.83"
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("D35").Value = "7.16"
$ws.Range("E35").Value = "  +11.27%  "
$ws.Range("D36").Value = "0.0622"
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "3.55"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "18.24"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "0.0231"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("D42").Value = "103.01"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("D43").Value = "1.536.13"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("E44").Value = "  +6.19%  "
$ws.Range("D45").Value = "1.11"
$ws.Range("E45").Value = "  +6.97%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.0921"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "7.83"
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").Value = "2.358.64"
$ws.Range("E50").Value = "  +3.15%  "
$ws.Range("D51").Value = "2.97"
$ws.Range("E51").Value = "  -0.06%  "
